# Auto-generated Excel COM-interop edit script
# Applies scheduled-runner profit recalculation updates across all 8 sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I2").Value = 484982.34
$ws.Range("K2").Value = 484982.34
$ws.Range("M2").Value = -484869.34
$ws.Range("H2").Value = 346886.94
$ws.Range("L5").Value = 733
$ws.Range("N5").Value = -963
$ws.Range("J5").Value = 733
$ws.Range("M5").Value = 15
$ws.Range("I5").Value = 100
$ws.Range("H5").Value = 574.75
$ws.Range("K5").Value = 100
$ws.Range("M18").Value = -340.8333
$ws.Range("K18").Value = 624.8333
$ws.Range("H18").Value = 624.8333
$ws.Range("I18").Value = 624.8333
$ws.Range("N69").Value = -47747.999
$ws.Range("L69").Value = 45999.999
$ws.Range("J69").Value = 15333.333
$ws.Range("H69").Value = 15333.333
$ws.Range("J72").Value = 15333.333
$ws.Range("H72").Value = 15333.333
$ws.Range("N72").Value = -146735.997
$ws.Range("L72").Value = 137999.997
$ws.Range("M138").Value = -1470.4634
$ws.Range("K138").Value = 6610.4634
$ws.Range("H138").Value = 2487.2104
$ws.Range("I138").Value = 2203.4878

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L5").Value = 400
$ws.Range("N5").Value = -624
$ws.Range("J5").Value = 400
$ws.Range("M5").Value = -51.40000000000001
$ws.Range("I5").Value = 163.4
$ws.Range("H5").Value = 301.41666
$ws.Range("K5").Value = 163.4
$ws.Range("M32").Value = -36071
$ws.Range("H32").Value = 32369.457
$ws.Range("K32").Value = 36358
$ws.Range("I32").Value = 36358
$ws.Range("N45").Value = -5353
$ws.Range("H45").Value = 4164.5
$ws.Range("J45").Value = 4599
$ws.Range("L45").Value = 4599
$ws.Range("L97").Value = 999
$ws.Range("H97").Value = 805.381
$ws.Range("J97").Value = 999
$ws.Range("N97").Value = -1991
$ws.Range("H122").Value = 6250
$ws.Range("K122").Value = 12000
$ws.Range("I122").Value = 4000
$ws.Range("M122").Value = -9550
$ws.Range("J132").Value = 1978.6923
$ws.Range("N132").Value = -10996.0769
$ws.Range("K132").Value = 62822.889
$ws.Range("L132").Value = 5936.0769
$ws.Range("M132").Value = -60292.889
$ws.Range("H132").Value = 17261.717
$ws.Range("I132").Value = 20940.963

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K4").Value = 163.4
$ws.Range("I4").Value = 163.4
$ws.Range("M4").Value = -48.40000000000001
$ws.Range("L4").Value = 400
$ws.Range("N4").Value = -630
$ws.Range("H4").Value = 301.41666
$ws.Range("J4").Value = 400
$ws.Range("L22").Value = 7907.8
$ws.Range("J22").Value = 7907.8
$ws.Range("N22").Value = -8253.799999999999
$ws.Range("H22").Value = 5061.5293
$ws.Range("L50").Value = 150000
$ws.Range("N50").Value = -151148
$ws.Range("H50").Value = 150000
$ws.Range("J50").Value = 150000
$ws.Range("H86").Value = 3691.8064
$ws.Range("I86").Value = 2535.5334
$ws.Range("M86").Value = -1412.5334
$ws.Range("K86").Value = 2535.5334
$ws.Range("K89").Value = 12677.667
$ws.Range("H89").Value = 3691.8064
$ws.Range("I89").Value = 2535.5334
$ws.Range("M89").Value = -7061.666999999999
$ws.Range("H134").Value = 4539
$ws.Range("L134").Value = 17232.666
$ws.Range("M134").Value = -9048.1875
$ws.Range("N134").Value = -22302.666
$ws.Range("I134").Value = 3861.0625
$ws.Range("K134").Value = 11583.1875
$ws.Range("J134").Value = 5744.222

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K22").Value = 349.25
$ws.Range("I22").Value = 349.25
$ws.Range("M22").Value = 0.75
$ws.Range("H22").Value = 994.3
$ws.Range("H28").Value = 9888
$ws.Range("N28").Value = -10378
$ws.Range("J28").Value = 9888
$ws.Range("L28").Value = 9888
$ws.Range("J41").Value = 9999.25
$ws.Range("L41").Value = 9999.25
$ws.Range("N41").Value = -10855.25
$ws.Range("H41").Value = 7666.1665
$ws.Range("L53").Value = 74999.5
$ws.Range("H53").Value = 74999.5
$ws.Range("N53").Value = -76213.5
$ws.Range("J53").Value = 74999.5
$ws.Range("M58").Value = -101987.6
$ws.Range("J58").Value = 3372.25
$ws.Range("L58").Value = 3372.25
$ws.Range("H58").Value = 73956.78999999999
$ws.Range("N58").Value = -3778.25
$ws.Range("I58").Value = 102190.6
$ws.Range("K58").Value = 102190.6
$ws.Range("N122").ClearContents()
$ws.Range("H122").Value = 2073.6667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6221.000100000001
$ws.Range("L122").Value = 0
$ws.Range("I122").Value = 2073.6667
$ws.Range("M122").Value = -3771.000100000001
$ws.Range("H136").Value = 73956.78999999999
$ws.Range("J136").Value = 3372.25
$ws.Range("L136").Value = 10116.75
$ws.Range("K136").Value = 306571.8
$ws.Range("M136").Value = -304021.8
$ws.Range("N136").Value = -15216.75
$ws.Range("I136").Value = 102190.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K6").Value = 29.25
$ws.Range("H6").Value = 9.75
$ws.Range("I6").Value = 9.75
$ws.Range("M6").Value = 83.75
$ws.Range("I11").Value = 278
$ws.Range("K11").Value = 834
$ws.Range("M11").Value = -694
$ws.Range("H11").Value = 652
$ws.Range("J105").Value = 9027.368
$ws.Range("N105").Value = -32324.104
$ws.Range("H105").Value = 9027.368
$ws.Range("L105").Value = 27082.104
$ws.Range("H122").Value = 1616.5
$ws.Range("K122").Value = 8995.5
$ws.Range("I122").Value = 999.5
$ws.Range("M122").Value = -6545.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5098.6665
$ws.Range("N70").Value = -5739.25
$ws.Range("L70").Value = 5199.25
$ws.Range("J70").Value = 5199.25
$ws.Range("K70").Value = 4897.5
$ws.Range("M70").Value = -4627.5
$ws.Range("I70").Value = 4897.5
$ws.Range("L73").Value = 5199.25
$ws.Range("H73").Value = 5098.6665
$ws.Range("K73").Value = 4897.5
$ws.Range("J73").Value = 5199.25
$ws.Range("M73").Value = -3961.5
$ws.Range("N73").Value = -7071.25
$ws.Range("I73").Value = 4897.5
$ws.Range("J80").Value = 21616.334
$ws.Range("L80").Value = 21616.334
$ws.Range("K80").Value = 0
$ws.Range("H80").Value = 21616.334
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -23612.334
$ws.Range("I80").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("L83").Value = 108081.67
$ws.Range("K83").Value = 0
$ws.Range("H83").Value = 21616.334
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -118065.67
$ws.Range("J83").Value = 21616.334
$ws.Range("K102").Value = 3690.7144
$ws.Range("M102").Value = -2068.7144
$ws.Range("N102").Value = -10243.75
$ws.Range("H102").Value = 4894
$ws.Range("I102").Value = 3690.7144
$ws.Range("L102").Value = 6999.75
$ws.Range("J102").Value = 6999.75
$ws.Range("N122").Value = -23641
$ws.Range("H122").Value = 2770
$ws.Range("J122").Value = 6247
$ws.Range("K122").Value = 7006.125
$ws.Range("L122").Value = 18741
$ws.Range("I122").Value = 2335.375
$ws.Range("M122").Value = -4556.125
$ws.Range("J132").Value = 2826
$ws.Range("N132").Value = -13538
$ws.Range("K132").Value = 155856.15
$ws.Range("L132").Value = 8478
$ws.Range("M132").Value = -153326.15
$ws.Range("H132").Value = 45544.305
$ws.Range("I132").Value = 51952.05

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2211.923
$ws.Range("K82").Value = 2031.2858
$ws.Range("M82").Value = -1670.2858
$ws.Range("N82").Value = -3144.6667
$ws.Range("J82").Value = 2422.6667
$ws.Range("I82").Value = 2031.2858
$ws.Range("L82").Value = 2422.6667
$ws.Range("L85").Value = 2422.6667
$ws.Range("I85").Value = 2031.2858
$ws.Range("J85").Value = 2422.6667
$ws.Range("N85").Value = -4918.6667
$ws.Range("M85").Value = -783.2858000000001
$ws.Range("K85").Value = 2031.2858
$ws.Range("H85").Value = 2211.923
$ws.Range("N122").Value = -19897
$ws.Range("H122").Value = 4495.0713
$ws.Range("J122").Value = 4999
$ws.Range("K122").Value = 9705.75
$ws.Range("L122").Value = 14997
$ws.Range("I122").Value = 3235.25
$ws.Range("M122").Value = -7255.75
$ws.Range("J132").Value = 6125.143
$ws.Range("N132").Value = -23435.429
$ws.Range("L132").Value = 18375.429
$ws.Range("H132").Value = 67069.16

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 58044.168
$ws.Range("N126").Value = -10749.9998
$ws.Range("I126").Value = 69265.664
$ws.Range("K126").Value = 207796.992
$ws.Range("L126").Value = 5809.9998
$ws.Range("M126").Value = -205326.992
$ws.Range("J126").Value = 1936.6666

Write-Host "Applied 221 value updates and 4 cell clears"